$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-5 are being rotated: the row data (columns D, K, L, M, N, O, P,
# Q, R, S, T) cyclically shifts: old row 2 -> new row 4, old row 4 -> new row 3,
# old row 3 -> new row 5, old row 5 -> new row 2. Columns A, B, C, E, F, G, H, I, J
# are identical across rows so they are unaffected.

# Capture the "before" values for the columns that actually change, for each row.
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

$rowData = @{}
foreach ($r in 2..5) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$r").Value2
    }
    $rowData[$r] = $vals
}

# Mapping of new row -> source old row
$mapping = @{
    2 = 5
    3 = 4
    4 = 2
    5 = 3
}

foreach ($newRow in 2..5) {
    $oldRow = $mapping[$newRow]
    $src = $rowData[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $src[$c]
    }
}
